$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Recompute A2 with a different pair of operands (0.1+0.7 instead of 0.6+0.3)
# so the classic floating-point rounding example still demonstrates the
# same 0.7999999999999993-ish binary rounding artifact.
$ws.Range("A2").Formula = "=0.1+0.7"

# Leave the active cell on A2 (matches where the edit was made).
$ws.Range("A2").Select()
